$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 34, pushing existing rows 34-35 down to 35-36.
$ws.Rows.Item(34).Insert()

# New row 34 gets a copy of row 33's ORIGINAL (pre-update) data.
$ws.Range("A34").Value = 1
$ws.Range("B34").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C34").Value = "Arica y Parinacota"
$ws.Range("D34").Value = 45008
$ws.Range("E34").Value = 15
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100103
$ws.Range("H34").Value = "Frutos de hueso (carozo)"
$ws.Range("I34").Value = 100103002
$ws.Range("J34").Value = "Ciruela"
$ws.Range("K34").Value = "Angeleno"
$ws.Range("L34").Value = "Segunda"
$ws.Range("M34").Value = 300
$ws.Range("N34").Value = 19000
$ws.Range("O34").Value = 20000
$ws.Range("P34").Value = 19500
$ws.Range("Q34").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R34").Value = "Región de O'Higgins"
$ws.Range("S34").Value = 1083
$ws.Range("T34").Value = 18

# Update row 33 (the original, most recent record) with new price data.
$ws.Range("D33").Value = 45013
$ws.Range("M33").Value = 280
$ws.Range("N33").Value = 21000
$ws.Range("O33").Value = 22000
$ws.Range("P33").Value = 21536
$ws.Range("S33").Value = 1196
